$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.03637433333333333
$ws.Range("H2").Value = 0.109123
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 1.908591666666666
$ws.Range("N2").Value = 5.725775
$ws.Range("O2").Value = 0.02908012806830959
$ws.Range("P2").Value = 0.0290801280683096
$ws.Range("Q2").Value = 0.06942374948055555
$ws.Range("R2").Value = 0.6248137453249999
$ws.Range("S2").Value = 0.02908012806830959
$ws.Range("T2").Value = 0.0290801280683096

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.03637433333333333
$ws.Range("H3").Value = 0.109123
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 20.524797
$ws.Range("N3").Value = 61.574391
$ws.Range("O3").Value = 0.3127246837341967
$ws.Range("P3").Value = 0.3127246837341967
$ws.Range("Q3").Value = 0.746575807677
$ws.Range("R3").Value = 6.719182269092999
$ws.Range("S3").Value = 0.3127246837341967
$ws.Range("T3").Value = 0.3127246837341967

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.03637433333333333
$ws.Range("H4").Value = 0.109123
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 8.719214000000001
$ws.Range("N4").Value = 26.157642
$ws.Range("O4").Value = 0.132849715422802
$ws.Range("P4").Value = 0.132849715422802
$ws.Range("Q4").Value = 0.3171555964406667
$ws.Range("R4").Value = 2.854400367966
$ws.Range("S4").Value = 0.132849715422802
$ws.Range("T4").Value = 0.132849715422802

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.03637433333333333
$ws.Range("H5").Value = 0.109123
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 11.88286
$ws.Range("N5").Value = 35.64858
$ws.Range("O5").Value = 0.1810523941044453
$ws.Range("P5").Value = 0.1810523941044453
$ws.Range("Q5").Value = 0.4322311105933334
$ws.Range("R5").Value = 3.89007999534
$ws.Range("S5").Value = 0.1810523941044453
$ws.Range("T5").Value = 0.1810523941044453

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.03637433333333333
$ws.Range("H6").Value = 0.109123
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 1.641831
$ws.Range("N6").Value = 4.925493
$ws.Range("O6").Value = 0.02501564718130951
$ws.Range("P6").Value = 0.02501564718130951
$ws.Range("Q6").Value = 0.059720508071
$ws.Range("R6").Value = 0.5374845726390001
$ws.Range("S6").Value = 0.02501564718130951
$ws.Range("T6").Value = 0.02501564718130951

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 0.03637433333333333
$ws.Range("H7").Value = 0.109123
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 20.954868
$ws.Range("N7").Value = 62.864604
$ws.Range("O7").Value = 0.3192774314889369
$ws.Range("P7").Value = 0.3192774314889369
$ws.Range("Q7").Value = 0.7622193535880001
$ws.Range("R7").Value = 6.859974182292
$ws.Range("S7").Value = 0.3192774314889369
$ws.Range("T7").Value = 0.3192774314889369
